{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that precedes it) that followed the\n// \"LOQ4084: Fen\u00f4menos de Transporte II (Requisito fraco)\" requirement line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"LOQ4084...\" requirement paragraph.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4084: Fen\u00f4menos de Transporte II\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // The three paragraphs immediately after the anchor are:\n  //   1) a blank paragraph\n  //   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n  //   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\"\n  // Collect them (guarding against a shorter document) and delete in\n  // reverse order so earlier indices stay valid while deleting.\n  const toDelete = [];\n  for (let offset = 1; offset <= 3 && anchorIndex + offset < items.length; offset++) {\n    toDelete.push(items[anchorIndex + offset]);\n  }\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph that precedes them) that followed the\n# \"LOQ4084: Fenomenos de Transporte II (Requisito fraco)\" requirement line.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*LOQ4084:*Transporte II*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1) {\n    # The three paragraphs immediately following the anchor are:\n    #   1) a blank paragraph\n    #   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n    #   3) \"(c) 2020 . Contact: luizeleno@usp.br. ...\"\n    # Delete all three in one shot by spanning from the end of the anchor\n    # paragraph to the end of the third paragraph after it.\n    $lastIndex = $anchorIndex + 3\n    if ($lastIndex -le $count) {\n        $startDel = $d.Paragraphs.Item($anchorIndex).Range.End\n        $endDel = $d.Paragraphs.Item($lastIndex).Range.End\n        $r = $d.Range($startDel, $endDel)\n        $r.Delete()\n    }\n}\n"}
